$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update failed test case indicators
$ws.Range("D17").Value = 79
$ws.Rows.Item(17).AutoFit()

$ws.Range("D21").Value = 5
$ws.Rows.Item(21).AutoFit()

# Update view / selection state (scroll position + active cell)
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H22").Select()
